$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F1 header text: "Raw Loudness Data" -> "Raw Pitch Data"
$ws.Range("F1").Value = "Raw Pitch Data"

# Update column F raw data values for the rows that changed
$ws.Range("F2").Value = 0.56479458917835679
$ws.Range("F3").Value = 0.60133266533066132
$ws.Range("F4").Value = 0.68515531062124246
$ws.Range("F5").Value = 0.51321142284569143
$ws.Range("F6").Value = 0.56049599198396793
$ws.Range("F7").Value = 0.4164929859719439
$ws.Range("F8").Value = 1.2418236472945892
$ws.Range("F9").Value = 0.64646793587174356
$ws.Range("F10").Value = 1.7039228456913826
$ws.Range("F11").Value = 0.64861723446893793
$ws.Range("F12").Value = 1.9446442885771544
$ws.Range("F13").Value = 0.63142284569138274
$ws.Range("F14").Value = 0.96456412825651305
$ws.Range("F15").Value = 0.82485971943887781
$ws.Range("F16").Value = 0.71954408817635274
$ws.Range("F17").Value = 1.2740631262525051
$ws.Range("F18").Value = 0.61637775551102203
$ws.Range("F19").Value = 0.78187374749498995
$ws.Range("F20").Value = 2.54
$ws.Range("F22").Value = 0.81315631262525057
$ws.Range("F23").Value = 0.80200400801603211
$ws.Range("F24").Value = 0.84661322645290582
$ws.Range("F25").Value = 0.87635270541082166
$ws.Range("F26").Value = 4.17
$ws.Range("F27").Value = 0.65330661322645289
$ws.Range("F28").Value = 4.17
$ws.Range("F29").Value = 1.0696593186372745
$ws.Range("F30").Value = 3.6020585615675791
$ws.Range("F31").Value = 0.84661322645290582
$ws.Range("F32").Value = 0.79085170340681354
$ws.Range("F33").Value = 0.71278557114228458
$ws.Range("F34").Value = 4.17
$ws.Range("F35").Value = 0.83174348697394784
$ws.Range("F36").Value = 0.76111222444889781
$ws.Range("F37").Value = 0.76854709418837686
$ws.Range("F38").Value = 0.96185370741482956
$ws.Range("F39").Value = 1.5975350701402804
$ws.Range("F40").Value = 4.17
$ws.Range("F42").Value = 1.2462725450901804
$ws.Range("F43").Value = 0.91679358717434867
$ws.Range("F44").Value = 1.2663627254509018
$ws.Range("F45").Value = 0.77214428857715434
$ws.Range("F46").Value = 0.77214428857715434
$ws.Range("F47").Value = 0.74803607214428858
$ws.Range("F48").Value = 4.75
$ws.Range("F49").Value = 1.4994088176352702
$ws.Range("F50").Value = 1.2181462925851705
$ws.Range("F51").Value = 1.1016232464929858
$ws.Range("F52").Value = 1.3306513026052105
$ws.Range("F53").Value = 1.2382364729458917
$ws.Range("F54").Value = 0.78018036072144281
$ws.Range("F55").Value = 1.0574248496993985
$ws.Range("F56").Value = 1.3708316633266533
$ws.Range("F57").Value = 1.3226152304609218
$ws.Range("F58").Value = 1.2060921843687371
$ws.Range("F59").Value = 1.0252805611222444
$ws.Range("F60").Value = 4.75
$ws.Range("F61").Value = 1.8610320641282563
$ws.Range("F62").Value = 1.606753507014028
$ws.Range("F63").Value = 1.216132264529058
$ws.Range("F64").Value = 1.198376753507014
$ws.Range("F65").Value = 0.89653306613226458
$ws.Range("F66").Value = 1.4291983967935871
$ws.Range("F67").Value = 0.87877755511022049
$ws.Range("F68").Value = 9.65
$ws.Range("F69").Value = 1.4558316633266533
$ws.Range("F70").Value = 6.1344088176352711
$ws.Range("F71").Value = 1.5978757515030058
$ws.Range("F72").Value = 1.7487975951903807
$ws.Range("F73").Value = 1.4380761523046091
$ws.Range("F74").Value = 1.3581763527054107
$ws.Range("F75").Value = 2.0506412825651301
$ws.Range("F76").Value = 1.1547950446347239
$ws.Range("F79").Value = 7.1553507014028064
$ws.Range("F80").Value = 1.5801202404809618
$ws.Range("F81").Value = 9.6056112224448906
$ws.Range("F82").Value = 6.59
$ws.Range("F83").Value = 6.59
$ws.Range("F84").Value = 0.87947895791583175
$ws.Range("F85").Value = 6.59
$ws.Range("F86").Value = 1.0310320641282567
$ws.Range("F87").Value = 0.55212424849699404
$ws.Range("F88").Value = 0.89160320641282576
$ws.Range("F89").Value = 1.0916533066132266
$ws.Range("F90").Value = 3.1588376753507017
$ws.Range("F91").Value = 1.1159018036072144
$ws.Range("F92").Value = 1.9403507014028054
$ws.Range("F93").Value = 6.59
$ws.Range("F94").Value = 6.59
$ws.Range("F95").Value = 1.6008717434869737
$ws.Range("F96").Value = 0.87341683366733469
$ws.Range("F97").Value = 3.183086172344689
$ws.Range("F98").Value = 1.2250200400801603
$ws.Range("F99").Value = 1.8191082164328656
$ws.Range("F100").Value = 1.443256513026052
$ws.Range("F102").Value = 1.3380961923847696
$ws.Range("F103").Value = 0.99390781563126251
$ws.Range("F104").Value = 1.0489779559118237
$ws.Range("F105").Value = 1.0902805611222446
$ws.Range("F106").Value = 1.0558617234468939
$ws.Range("F107").Value = 7.74
$ws.Range("F108").Value = 1.062745490981964
$ws.Range("F109").Value = 1.0076753507014029
$ws.Range("F110").Value = 1.0076753507014029
$ws.Range("F111").Value = 1.9232164328657317
$ws.Range("F112").Value = 1.3380961923847696
$ws.Range("F113").Value = 1.6065631262525053
$ws.Range("F114").Value = 1.207304609218437
$ws.Range("F115").Value = 1.1522344689378758
$ws.Range("F116").Value = 1.1522344689378758
$ws.Range("F117").Value = 1.110931863727455
$ws.Range("F118").Value = 1.8612625250501003
$ws.Range("F119").Value = 1.5721442885771544
$ws.Range("F120").Value = 1.3243286573146293
$ws.Range("F121").Value = 6.4966194889779567

# Update sheet view: select F1 (clears topLeftCell scroll position)
$ws.Range("F1").Select()
